# "Tareas diarias": row 5 (Visualizacion productos, seccion TIENDA) is now
# finished, so restyle it to match the "Completo" rows (same look as row 6)
# and flip its ESTADO column from "Incompleto" to "Completo". Then move the
# active selection down to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tareas diarias")

# Row 6 already carries the "done" formatting (green "Buena" cell style per
# column) we want on row 5 - copy just the formats over, cell by cell, so
# each column keeps its own alignment/wrap behaviour.
$ws.Range("A6").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B6").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Mark the task as complete.
$ws.Range("C5").Value = "Completo"

# Update the saved selection/view state.
$ws.Range("D7").Select()
